{"js": "// The diff rotates 8 blocks of text across 8 paragraph \"slots\" in the\n// document. Because we know the exact old/new text for every slot, the\n// simplest, most robust approach is a direct, paragraph-scoped\n// find & replace for each slot (no generic \"rotation\" logic needed).\n//\n// Word manual line breaks (<w:br/>) are represented in Office.js range\n// text as the vertical-tab character \"\\v\" (U+000B); inserting text that\n// contains \"\\v\" produces <w:t>...</w:t><w:br/><w:t>...</w:t> inside a\n// single run, matching the target OOXML exactly. A literal \"\\n\" would\n// instead split into a new paragraph, so it is never used here.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Replace the entire text of a whole paragraph (single-run paragraphs).\nfunction replaceParagraph(paragraph, newText) {\n  paragraph.insertText(newText, Word.InsertLocation.replace);\n}\n\n// Replace one specific old snippet of text inside a paragraph that\n// contains multiple runs (e.g. bold \"Label: \" runs followed by plain\n// value runs), leaving the rest of the paragraph untouched. Requires the\n// snippet to currently be unique within the paragraph.\nasync function replaceWithin(paragraph, oldText, newText) {\n  const found = paragraph.search(oldText, { matchCase: true });\n  found.load(\"items\");\n  await context.sync();\n  if (found.items.length !== 1) {\n    throw new Error(\n      \"replaceWithin: expected exactly 1 match for \" +\n        JSON.stringify(oldText.substring(0, 40)) +\n        \", got \" +\n        found.items.length\n    );\n  }\n  found.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- Slot 1: \"Objetivos\" section body ---------------------------------\nreplaceParagraph(\n  paragraphs.items[5],\n  \"Desenvolvimento do trabalho de conclus\u00e3o de curso, sob orienta\u00e7\u00e3o de um professor orientador,\\v\" +\n    \"o qual deve constituir-se num projeto de tema espec\u00edfico relacionado \u00e0s atribui\u00e7\u00f5es da profiss\u00e3o.\"\n);\n\n// --- Slot 2: \"Docente(s) Respons\u00e1vel(eis)\" entry -----------------------\nreplaceParagraph(\n  paragraphs.items[7],\n  \"Possibilitar aos alunos a realiza\u00e7\u00e3o de trabalho de s\u00edntese e integra\u00e7\u00e3o dos conhecimentos adquiridos ao longo do curso, conforme projeto aprovado na disciplina de Trabalho de Conclus\u00e3o do Curso I.\"\n);\n\n// --- Slot 3: \"Programa resumido\" body -----------------------------------\nreplaceParagraph(\n  paragraphs.items[9],\n  \"Elabora\u00e7\u00e3o de uma monografia ou de relat\u00f3rio t\u00e9cnico que apresente: (1) o tema e sua\\v\" +\n    \"import\u00e2ncia, (2) os objetivos, (3) a revis\u00e3o bibliogr\u00e1fica, (4) a metodologia cient\u00edfica (5) o\\v\" +\n    \"desenvolvimento do projeto, (6) a an\u00e1lise e discuss\u00e3o dos resultados, (7) as conclus\u00f5es e (8)\\v\" +\n    \"refer\u00eancias bibliogr\u00e1ficas\"\n);\n\n// --- Slot 4: \"Programa\" body --------------------------------------------\nreplaceParagraph(\n  paragraphs.items[11],\n  \"Reuni\u00f5es peri\u00f3dicas com o orientador e realiza\u00e7\u00e3o do trabalho conforme orienta\u00e7\u00e3o e\\v\" +\n    \"apresenta\u00e7\u00e3o de uma monografia final, conforme norma do Departamento de Engenharia Qu\u00edmica\\v\" +\n    \"e Produ\u00e7\u00e3o.\"\n);\n\nawait context.sync();\n\n// --- \"Avalia\u00e7\u00e3o\" paragraph holds 3 labeled values in one paragraph -----\n// New text for \"M\u00e9todo\" happens to equal the old \"Crit\u00e9rio\" text (and\n// likewise Crit\u00e9rio's new text equals Norma's old text), so the three\n// replacements MUST run back-to-front; otherwise an earlier replacement\n// would make its own snippet ambiguous (no longer unique) for a later\n// search still looking for the original wording.\nconst avaliacao = paragraphs.items[13];\n\n// Slot 7: Norma de recupera\u00e7\u00e3o value\nawait replaceWithin(\n  avaliacao,\n  \"Reavalia\u00e7\u00e3o, preferencialmente para a mesma banca, com as modifica\u00e7\u00f5es sugeridas.\",\n  \"Recomendada pelo orientador\"\n);\n\n// Slot 6: Crit\u00e9rio value\nawait replaceWithin(\n  avaliacao,\n  \"Avalia\u00e7\u00e3o perante uma banca examinadora composta por 3 (tr\u00eas) membros, conforme norma do\\v\" +\n    \"Departamento de Engenharia Qu\u00edmica e Produ\u00e7\u00e3o.\",\n  \"Reavalia\u00e7\u00e3o, preferencialmente para a mesma banca, com as modifica\u00e7\u00f5es sugeridas.\"\n);\n\n// Slot 5: M\u00e9todo value\nawait replaceWithin(\n  avaliacao,\n  \"Reuni\u00f5es peri\u00f3dicas com o orientador e realiza\u00e7\u00e3o do trabalho conforme orienta\u00e7\u00e3o e\\v\" +\n    \"apresenta\u00e7\u00e3o de uma monografia final, conforme norma do Departamento de Engenharia Qu\u00edmica\\v\" +\n    \"e Produ\u00e7\u00e3o.\",\n  \"Avalia\u00e7\u00e3o perante uma banca examinadora composta por 3 (tr\u00eas) membros, conforme norma do\\v\" +\n    \"Departamento de Engenharia Qu\u00edmica e Produ\u00e7\u00e3o.\"\n);\n\n// --- Slot 8: \"Bibliografia\" body ----------------------------------------\nreplaceParagraph(paragraphs.items[15], \"1285870 - Marcos Villela Barcza\");\n\nawait context.sync();\n", "ps1": "# The diff rotates 8 blocks of text across 8 paragraph \"slots\" in the\n# document. Because we know the exact old/new text for every slot, the\n# simplest, most robust approach is a direct, paragraph-scoped\n# find & replace for each slot (no generic \"rotation\" logic needed).\n#\n# Word manual line breaks (<w:br/>) are represented in Range text as the\n# vertical-tab character (PowerShell escape \"`v\", U+000B); assigning /\n# replacing with text that contains \"`v\" produces\n# <w:t>...</w:t><w:br/><w:t>...</w:t> inside a single run, matching the\n# target OOXML exactly.\n#\n# Note: COM Range objects are re-fetched fresh (via $d.Paragraphs.Item(n).Range)\n# before every Find/Replace call instead of being reused or passed through\n# helper functions, because a Range narrows to the just-found text after a\n# successful Find, and passing a live Range through a PowerShell function\n# parameter does not marshal correctly in this host.\n\n$d = $word.ActiveDocument\n\n# Word \"Replace\" constant: wdReplaceOne = 1\n$wdReplaceOne = 1\n\n# --- Slot 1: \"Objetivos\" section body (paragraph 6) ---------------------\n$d.Paragraphs.Item(6).Range.Text = \"Desenvolvimento do trabalho de conclus\u00e3o de curso, sob orienta\u00e7\u00e3o de um professor orientador,`v\" + `\n    \"o qual deve constituir-se num projeto de tema espec\u00edfico relacionado \u00e0s atribui\u00e7\u00f5es da profiss\u00e3o.\"\n\n# --- Slot 2: \"Docente(s) Respons\u00e1vel(eis)\" entry (paragraph 8) ----------\n$d.Paragraphs.Item(8).Range.Text = \"Possibilitar aos alunos a realiza\u00e7\u00e3o de trabalho de s\u00edntese e integra\u00e7\u00e3o dos conhecimentos adquiridos ao longo do curso, conforme projeto aprovado na disciplina de Trabalho de Conclus\u00e3o do Curso I.\"\n\n# --- Slot 3: \"Programa resumido\" body (paragraph 10) ---------------------\n$d.Paragraphs.Item(10).Range.Text = \"Elabora\u00e7\u00e3o de uma monografia ou de relat\u00f3rio t\u00e9cnico que apresente: (1) o tema e sua`v\" + `\n    \"import\u00e2ncia, (2) os objetivos, (3) a revis\u00e3o bibliogr\u00e1fica, (4) a metodologia cient\u00edfica (5) o`v\" + `\n    \"desenvolvimento do projeto, (6) a an\u00e1lise e discuss\u00e3o dos resultados, (7) as conclus\u00f5es e (8)`v\" + `\n    \"refer\u00eancias bibliogr\u00e1ficas\"\n\n# --- Slot 4: \"Programa\" body (paragraph 12) ------------------------------\n$d.Paragraphs.Item(12).Range.Text = \"Reuni\u00f5es peri\u00f3dicas com o orientador e realiza\u00e7\u00e3o do trabalho conforme orienta\u00e7\u00e3o e`v\" + `\n    \"apresenta\u00e7\u00e3o de uma monografia final, conforme norma do Departamento de Engenharia Qu\u00edmica`v\" + `\n    \"e Produ\u00e7\u00e3o.\"\n\n# --- \"Avalia\u00e7\u00e3o\" paragraph (14) holds 3 labeled values in one paragraph -\n# New text for \"M\u00e9todo\" happens to equal the old \"Crit\u00e9rio\" text (and\n# likewise Crit\u00e9rio's new text equals Norma's old text), so the three\n# replacements MUST run back-to-front; otherwise an earlier replacement\n# would make its own snippet ambiguous for a later search still looking\n# for the original wording.\n\n# Slot 7: Norma de recupera\u00e7\u00e3o value\n$rNorma = $d.Paragraphs.Item(14).Range\n$rNorma.Find.Execute(\n    \"Reavalia\u00e7\u00e3o, preferencialmente para a mesma banca, com as modifica\u00e7\u00f5es sugeridas.\",\n    $true, $true, $false, $false, $false, $true, 0, $false,\n    \"Recomendada pelo orientador\", $wdReplaceOne) | Out-Null\n\n# Slot 6: Crit\u00e9rio value\n$rCriterio = $d.Paragraphs.Item(14).Range\n$rCriterio.Find.Execute(\n    (\"Avalia\u00e7\u00e3o perante uma banca examinadora composta por 3 (tr\u00eas) membros, conforme norma do`v\" + `\n        \"Departamento de Engenharia Qu\u00edmica e Produ\u00e7\u00e3o.\"),\n    $true, $true, $false, $false, $false, $true, 0, $false,\n    \"Reavalia\u00e7\u00e3o, preferencialmente para a mesma banca, com as modifica\u00e7\u00f5es sugeridas.\", $wdReplaceOne) | Out-Null\n\n# Slot 5: M\u00e9todo value\n$rMetodo = $d.Paragraphs.Item(14).Range\n$rMetodo.Find.Execute(\n    (\"Reuni\u00f5es peri\u00f3dicas com o orientador e realiza\u00e7\u00e3o do trabalho conforme orienta\u00e7\u00e3o e`v\" + `\n        \"apresenta\u00e7\u00e3o de uma monografia final, conforme norma do Departamento de Engenharia Qu\u00edmica`v\" + `\n        \"e Produ\u00e7\u00e3o.\"),\n    $true, $true, $false, $false, $false, $true, 0, $false,\n    (\"Avalia\u00e7\u00e3o perante uma banca examinadora composta por 3 (tr\u00eas) membros, conforme norma do`v\" + `\n        \"Departamento de Engenharia Qu\u00edmica e Produ\u00e7\u00e3o.\"), $wdReplaceOne) | Out-Null\n\n# --- Slot 8: \"Bibliografia\" body (paragraph 16) --------------------------\n$d.Paragraphs.Item(16).Range.Text = \"1285870 - Marcos Villela Barcza\"\n"}
